$d = $word.ActiveDocument

# Locate the "Make a todo list project" bullet (the last populated item
# in the numbered list before the trailing blank bullet) and use it as
# the anchor to insert the five new sub-items about passing data from a
# parent component to a child component.
$anchor = $d.Content
$anchor.Find.Execute("list project", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$anchor.Expand(4) | Out-Null

# wdListLevelNumber values are 1-based: 1 == top-level (ilvl 0 in OOXML),
# 2 == first nested level (ilvl 1 in OOXML).

# 1) Top-level bullet: "Send data from parent component to child component"
$anchor.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p1.Range.InsertAfter("Send data from parent component to child component")

# 2) Nested bullet: "Make a child component"
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p2.Range.InsertAfter("Make a child component")
$p2.Range.ListFormat.ListLevelNumber = 2

# 3) Nested bullet: "Use the child component in the parent component "
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p3.Range.InsertAfter("Use the child component in the parent component ")
$p3.Range.ListFormat.ListLevelNumber = 2

# 4) Nested bullet: "Pass data from parent to child component "
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p4.Range.InsertAfter("Pass data from parent to child component ")
$p4.Range.ListFormat.ListLevelNumber = 2

# 5) Nested bullet: "Import input decorator in child component and capture the data."
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p5.Range.InsertAfter("Import input decorator in child component and capture the data.")
$p5.Range.ListFormat.ListLevelNumber = 2
